# #5: property boat&car done
# Insert a new "建物" (building) worksheet between "土地" (land) and "債務"
# (debt), populated with one property record, mirroring the layout of the
# "土地" sheet (header row in B1:Q1, data starting row 2).

$wb = $excel.ActiveWorkbook

$landSheet = $wb.Worksheets.Item("土地")

# Create the new worksheet right after "土地" so the final tab order is
# 土地, 建物, 債務, 事業投資.
$buildingSheet = $wb.Worksheets.Add($null, $landSheet)
$buildingSheet.Name = "建物"

# Copy header row (with its bold/bordered style) and the first data row
# (with its style) from 土地 as a formatting template, then overwrite the
# values for the new building record.
$landSheet.Range("B1:Q1").Copy($buildingSheet.Range("B1:Q1"))
$landSheet.Range("A2:Q2").Copy($buildingSheet.Range("A2:Q2"))

$buildingSheet.Range("A2").Value = 22
$buildingSheet.Range("B2").Value = "臺中市沙鹿區屏西路"
$buildingSheet.Range("C2").Value = 432
$buildingSheet.Range("D2").Value = "全部"
$buildingSheet.Range("E2").Value = "顔清標"

# "83年10月14H" is plain text (a ROC-calendar-style date, not a recognised
# Excel date) so it is safe to assign directly.
$buildingSheet.Range("F2").Value = "83年10月14H"

$buildingSheet.Range("G2").Value = "買賣"
$buildingSheet.Range("H2").Value = "(超過五年）"
$buildingSheet.Range("I2").Value = "land"
$buildingSheet.Range("J2").Value = "normal"

# "2012-04-10" looks like an ISO date, so Excel would otherwise silently
# reinterpret it as a date serial number. Force the cell to text first so
# it is stored verbatim, matching the other sheets where this same string
# is kept as a shared text string.
$buildingSheet.Range("K2").NumberFormat = "@"
$buildingSheet.Range("K2").Value = "2012-04-10"

$buildingSheet.Range("L2").Value = "顏清標"
$buildingSheet.Range("M2").Value = 979
$buildingSheet.Range("N2").Value = "tmp1b4d1"
$buildingSheet.Range("O2").Value = 22
$buildingSheet.Range("P2").Value = 1
$buildingSheet.Range("Q2").Value = 432

# Keep "土地" as the selected/active tab (it was active before this edit;
# adding a sheet would otherwise switch selection to the new sheet).
$landSheet.Activate()
